$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 10 - "Chart" component definition (dropdown comp refactor + chart comp generation)
$ws.Range("A10").Value = "Chart"
$ws.Range("C10").Value = "type"
$ws.Range("D10").Value = "Bar"
$ws.Range("E10").Value = "label"
$ws.Range("F10").Value = "Renewed Policies"
$ws.Range("G10").Value = "data"
$ws.Range("I10").Value = "options"
$ws.Range("J10").Value = '{"maintainAspectRatio": false, "scales":{"yAxes":[{"ticks":{"beginAtZero":true}}]}}'
$ws.Range("H10").Value = '{"labels":["January", "February", "March", "April", "May", "June"],"datasets":[{"label":"# of Policies Renewed","data":[12,19,3,5,2,3],"backgroundColor":["rgba(255, 99, 132, 0.2)","rgba(54, 162, 235, 0.2)","rgba(255, 206, 86, 0.2)","rgba(75, 192, 192, 0.2)","rgba(153, 102, 255, 0.2)","rgba(255, 159, 64, 0.2)"],"borderColor":["rgba(255, 99, 132, 1)","rgba(54, 162, 235, 1)","rgba(255, 206, 86, 1)","rgba(75, 192, 192, 1)","rgba(153, 102, 255, 1)","rgba(255, 159, 64, 1)"],"borderWidth":1}]}'

# Cells H10 / I10 carry an explicit "General" number format in the source workbook
$ws.Range("H10").NumberFormat = "General"
$ws.Range("I10").NumberFormat = "General"

# Selection moved to the newly added cell
$ws.Range("H10").Select()
